$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be forced back to
# Text format first, otherwise Excel auto-converts the assigned string into a
# numeric value (losing the original inline-string / Text cell type).
# Row 2
$ws.Range("D2").Value = "26.618.54"
$ws.Range("E2").Value = "  -0.13%  "
# Row 3
$ws.Range("D3").Value = "1.596.60"
$ws.Range("E3").Value = "  +0.09%  "
# Row 4
$ws.Range("E4").Value = "  +0.25%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.04"
$ws.Range("E5").Value = "  -0.44%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.514"
$ws.Range("E6").Value = "  +0.13%  "
# Row 7
$ws.Range("E7").Value = "  +0.23%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0617"
$ws.Range("E8").Value = "  -0.02%  "
# Row 9
$ws.Range("E9").Value = "  -0.36%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.45"
$ws.Range("E10").Value = "  -1.26%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0836"
$ws.Range("E11").Value = "  +0.05%  "
# Row 12
$ws.Range("D12").Value = "1.820.99"
$ws.Range("E12").Value = "  +0.18%  "
# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.02"
$ws.Range("E13").Value = "  -0.16%  "
# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.561.21"
$ws.Range("E14").Value = "  -2.10%  "
# Row 15
$ws.Range("E15").Value = "  -0.43%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.99"
$ws.Range("E16").Value = "  -0.30%  "
# Row 17
$ws.Range("D17").Value = "26.614.56"
$ws.Range("E17").Value = "  +0.01%  "
# Row 18
$ws.Range("E18").Value = "  +0.71%  "
# Row 19
$ws.Range("E19").Value = "  +0.17%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.32"
$ws.Range("E20").Value = "  -1.09%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.02"
$ws.Range("E21").Value = "  +4.69%  "
# Row 22
$ws.Range("E22").Value = "  +0.28%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.29"
$ws.Range("E23").Value = "  -0.66%  "
# Row 24
$ws.Range("E24").Value = "  +0.14%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.18"
$ws.Range("E25").Value = "  -0.85%  "
# Row 26
$ws.Range("E26").Value = "  +0.14%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.13"
$ws.Range("E27").Value = "  -0.40%  "
# Row 28
$ws.Range("E28").Value = "  +0.00%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.24"
$ws.Range("E29").Value = "  -0.71%  "
# Row 30
$ws.Range("E30").Value = "  +0.20%  "
# Row 31
$ws.Range("E31").Value = "  +0.06%  "
# Row 32
$ws.Range("E32").Value = "  -0.53%  "
# Row 33
$ws.Range("E33").Value = "  +0.64%  "
# Row 34
$ws.Range("D34").Value = "1.274.91"
$ws.Range("E34").Value = "  -1.56%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.615"
$ws.Range("E35").Value = "  -8.18%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.46"
$ws.Range("E36").Value = "  +0.33%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.48"
$ws.Range("E37").Value = "  -0.19%  "
# Row 38
$ws.Range("E38").Value = "  -0.95%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.837"
$ws.Range("E39").Value = "  +0.16%  "
# Row 40
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.978"
$ws.Range("E40").Value = "  +16.82%  "
# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.47"
$ws.Range("E41").Value = "  +1.65%  "
# Row 42
$ws.Range("E42").Value = "  +0.18%  "
# Row 43
$ws.Range("E43").Value = "  -0.83%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.94"
$ws.Range("E44").Value = "  +0.27%  "
# Row 45
$ws.Range("D45").Value = "1.733.35"
$ws.Range("E45").Value = "  +0.19%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.02"
$ws.Range("E46").Value = "  +0.37%  "
# Row 47
$ws.Range("E47").Value = "  -0.06%  "
# Row 48
$ws.Range("E48").Value = "  +3.52%  "
# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.46"
$ws.Range("E50").Value = "  -0.50%  "
# Row 51
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.14%  "
